# Actualización automática de tasas-transfi.xlsx

$wb = $excel.ActiveWorkbook

# --- Update text on "Hoja1" (A1): refreshed conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.1 = 28360.54 pesos`n✅ 28360.54 pesos = 7.03 = 959.88 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update numeric rate cells on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 140.9
$wsTasas.Range("O10").Value = 3996
$wsTasas.Range("N12").Value = 4033
$wsTasas.Range("O12").Value = 136.5
